$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''246.00'
$ws.Range("E2").Value = '''-0.62%'
$ws.Range("D3").Value = '''29.81'
$ws.Range("E3").Value = '''-1.53%'
$ws.Range("D4").Value = '''5.152'
$ws.Range("E4").Value = '''-0.52%'
$ws.Range("D5").Value = '''0.05770'
$ws.Range("E5").Value = '''0.35%'
$ws.Range("E6").Value = '''1.00%'
$ws.Range("D7").Value = '''3.239'
$ws.Range("E7").Value = '''6.78%'
$ws.Range("D8").Value = '''0.8501'
$ws.Range("E8").Value = '''-1.28%'
$ws.Range("D9").Value = '''0.8531'
$ws.Range("E9").Value = '''-2.40%'
$ws.Range("D10").Value = '''0.1380'
$ws.Range("E10").Value = '''1.04%'
$ws.Range("D11").Value = '''0.07083'
$ws.Range("E11").Value = '''1.34%'
$ws.Range("D12").Value = '''0.03258'
$ws.Range("E12").Value = '''11.72%'
$ws.Range("D13").Value = '''0.09368'
$ws.Range("E13").Value = '''-0.29%'
$ws.Range("D14").Value = '''0.001537'
$ws.Range("E14").Value = '''1.90%'
$ws.Range("D15").Value = '''0.0005948'
$ws.Range("E15").Value = '''-94.17%'
$ws.Range("D16").Value = '''0.006005'
$ws.Range("E16").Value = '''-2.27%'
$ws.Range("D17").Value = '''3.515'
$ws.Range("E17").Value = '''0.27%'
$ws.Range("D18").Value = '''2.223'
$ws.Range("E18").Value = '''-2.30%'
$ws.Range("D19").Value = '''0.3163'
$ws.Range("E19").Value = '''-0.70%'
$ws.Range("E20").Value = '''2.05%'
$ws.Range("D21").Value = '''0.1316'
$ws.Range("E21").Value = '''0.63%'
$ws.Range("D22").Value = '''3.481'
$ws.Range("E22").Value = '''-3.52%'
$ws.Range("B23").Value = 'CoinExToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D23").Value = '''0.04135'
$ws.Range("E23").Value = '''0.24%'
$ws.Range("B24").Value = 'ZBToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D24").Value = '''0.1410'
$ws.Range("E24").Value = '''2.34%'
$ws.Range("D25").Value = '''0.001227'
$ws.Range("E25").Value = '''1.17%'
$ws.Range("E26").Value = '''-8.12%'
$ws.Range("E27").Value = '''1.84%'
$ws.Range("D28").Value = '''0.0001448'
$ws.Range("E28").Value = '''4.20%'
$ws.Range("D40").Value = '''0.03745'
$ws.Range("E40").Value = '''-1.14%'
$ws.Range("B41").Value = 'BKEXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D41").Value = '''0.1071'
$ws.Range("E41").Value = '''0.04%'
$ws.Range("B42").Value = 'CEJI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D42").Value = '''0.002470'
$ws.Range("E42").Value = '''7.53%'
$ws.Range("B43").Value = 'KickToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D43").Value = '''0.002950'
$ws.Range("E43").Value = '''-48.09%'
$ws.Range("D44").Value = '''0.008477'
$ws.Range("E44").Value = '''-16.58%'
$ws.Range("D45").Value = '''0.00005494'
$ws.Range("E45").Value = '''7.52%'
$ws.Range("E46").Value = '''0.11%'
$ws.Range("D47").Value = '''0.07097'
$ws.Range("E47").Value = '''-20.19%'
$ws.Range("D48").Value = '''0.002232'
$ws.Range("E48").Value = '''-17.85%'
$ws.Range("D49").Value = '''0.00002099'
$ws.Range("E49").Value = '''0.11%'
$ws.Range("D50").Value = '''0.0001999'
$ws.Range("E50").Value = '''0.11%'
